# Update CodeSystem-identifier-type.xlsx metadata per new release.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 6.0.0 -> 6.1.0
$meta.Range("B3").Value = "6.1.0"

# Date: 2022-01-21T20:46:54+00:00 -> 2022-05-31T20:10:14+00:00
$meta.Range("B8").Value = "2022-05-31T20:10:14+00:00"

# Description: drop "IBM " before "Health Data Connect"
$meta.Range("B11").Value = "Extended set of Identifier type code for Health Data Connect Common Data Model resources"

# --- Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Definition for DRI row: "internal IBM process" -> "internal HDC process"
$concepts.Range("D2").Value = "Business identifier for derived resources created by an internal HDC process"
